$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Status (D) and Responsible Unit (C) columns for rows 12-21
# Rows 12-14, 16-19: Status "in progress" -> "done"
$ws.Range("D12").Value = "done"
$ws.Range("D13").Value = "done"
$ws.Range("D14").Value = "done"
$ws.Range("D16").Value = "done"
$ws.Range("D17").Value = "done"
$ws.Range("D18").Value = "done"
$ws.Range("D19").Value = "done"

# Row 15: was blank Responsible Unit/Status -> reject/reject
$ws.Range("C15").Value = "reject"
$ws.Range("D15").Value = "reject"

# Row 19: Responsible Unit blank -> lgoychev
$ws.Range("C19").Value = "lgoychev"

# Row 20: Responsible Unit blank -> all; Status blank -> done
$ws.Range("C20").Value = "all"
$ws.Range("D20").Value = "done"

# Row 21: Responsible Unit blank -> lgoychev; Status blank -> done
$ws.Range("C21").Value = "lgoychev"
$ws.Range("D21").Value = "done"

# Apply center horizontal alignment to the newly-filled/changed status & responsible-unit cells
$ws.Range("C15:D15").HorizontalAlignment = -4108
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("D18").HorizontalAlignment = -4108
$ws.Range("C19").HorizontalAlignment = -4108
$ws.Range("D19").HorizontalAlignment = -4108
$ws.Range("C20:D20").HorizontalAlignment = -4108
$ws.Range("C21:D21").HorizontalAlignment = -4108

# Move cell selection to E24 to match end state cursor position
$ws.Range("E24").Select()
